$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New PSSM score matrix values (B2:K21), row-major, columns B..K
$values = @(
    ,(-17.34350671825161, -0.4066272768562257, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161)
    ,(-17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, 0.8225435513393097, -17.34350671825161, -17.34350671825161)
    ,(-17.34350671825161, -0.1602881255799057, 0.5268881053312987, -17.34350671825161, 3.880832412009773, -17.34350671825161, 1.286621734134304, -17.34350671825161, -17.34350671825161, -17.34350671825161)
    ,(-17.34350671825161, 0.3676560950042475, -17.34350671825161, -17.34350671825161, -17.34350671825161, 3.439214005409011, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161)
    ,(-17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161)
    ,(3.229997170333371, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161)
    ,(-17.34350671825161, -17.34350671825161, -17.34350671825161, 1.46488248060533, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161)
    ,(3.40833559416175, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161)
    ,(-17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, 0.2836352621583946, -17.34350671825161, 1.95615073053282)
    ,(-17.34350671825161, -17.34350671825161, -17.34350671825161, 1.939963768809139, -17.34350671825161, 1.420396143456054, -17.34350671825161, -17.34350671825161, -17.34350671825161, 1.984705811381224)
    ,(-17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161)
    ,(-17.34350671825161, -17.34350671825161, -17.34350671825161, 2.010641487574297, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, 1.805021997820626)
    ,(-17.34350671825161, -17.34350671825161, 0.9970605924212776, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, 2.51205680962203)
    ,(-17.34350671825161, -17.34350671825161, -0.3867562312949653, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161)
    ,(-17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161)
    ,(-17.34350671825161, 1.072008125887715, -0.02104158711604699, -17.34350671825161, -17.34350671825161, -17.34350671825161, 1.537681942603843, 1.788041310048033, 4.321919853816804, -17.34350671825161)
    ,(-17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, -17.34350671825161, 2.295411705697136, 0.01991797638700582, -17.34350671825161, -17.34350671825161)
    ,(-17.34350671825161, -17.34350671825161, 2.990675985439187, -17.34350671825161, -17.34350671825161, -17.34350671825161, 2.099430704817633, 1.281563528167942, -17.34350671825161, -17.34350671825161)
    ,(-17.34350671825161, 3.097322919038444, 2.779120354746611, -17.34350671825161, 2.397362686318584, -17.34350671825161, 1.503884070572508, 3.338532011568755, -17.34350671825161, 1.567138939761803)
    ,(-17.34350671825161, 2.678022295325321, -17.34350671825161, 3.228498613208306, -17.34350671825161, 2.695188600401428, 1.393242023891653, -17.34350671825161, -17.34350671825161, -17.34350671825161)
)

$cols = @("B","C","D","E","F","G","H","I","J","K")
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = 2 + $i
    $rowValues = $values[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range($cols[$j] + $row).Value = $rowValues[$j]
    }
}
